$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find and delete the "Blood Hunter" row entirely (class name + description)
$ws.Rows.Item(4).Delete()

# Update the selection to match the post-edit state
$ws.Range("A4:B4").Select()
